$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint2")

# Enter the new retrospective content in the same order the author typed it
# (this keeps shared-string allocation order faithful to the original edit).
$ws.Range("A20").Value = "Make 2 lists:"
$ws.Range("A22").Value = "1. Reuse code "
$ws.Range("A23").Value = "2. Refactor code after each sprint"
$ws.Range("A24").Value = "3. Follow function naming standards"
$ws.Range("E21").Value = "Things that want to avoid doing in the future"
$ws.Range("A21").Value = "Things want to keep doing"
$ws.Range("A25").Value = "4. Write valid testcases"
$ws.Range("E22").Value = "1. Wait untill last moment"
$ws.Range("A26").Value = "5. Run testcases after integrating code in each sprint "
$ws.Range("E23").Value = "2. Avoid duplicate code"
$ws.Range("E24").Value = "3. Avoid complex code "
$ws.Range("E25").Value = "4. Avoid breaking testcases"
$ws.Range("A27").Value = "6. Update testdata for each usecase"

# Section title and the two list headings are bold, like the other headings
# used throughout this workbook.
$ws.Range("A20").Font.Bold = $true
$ws.Range("A21").Font.Bold = $true
$ws.Range("E21").Font.Bold = $true

# Small blank spacer cell above the new section, formatted with wrap text.
$ws.Range("A17").WrapText = $true

# Make Sprint2 the active/selected sheet and leave the selection where the
# author left it, which also drives workbook.xml's activeTab bookkeeping.
$ws.Activate()
$ws.Range("G29").Select()
